$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.988.48"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.911.54"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8134"
$ws.Range("E5").Value = "  +8.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.44"
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9987"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3140"
$ws.Range("E8").Value = "  +3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.57"
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06998"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08015"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7471"
$ws.Range("E12").Value = "  +0.62%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.906.46"
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.191"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.61"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.972.00"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.04"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.889"
$ws.Range("E18").Value = "  -0.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.58"
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007784"
$ws.Range("E20").Value = "  +1.57%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.152.06"
$ws.Range("E22").Value = "  +0.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9985"
$ws.Range("E23").Value = "  -0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.960"
$ws.Range("E24").Value = "  +0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1619"
$ws.Range("E25").Value = "  +26.93%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "168.47"
$ws.Range("E26").Value = "  +1.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.231"
$ws.Range("E27").Value = "  +0.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.89"
$ws.Range("E28").Value = "  +1.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.074"
$ws.Range("E29").Value = "  +2.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.362"
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.511"
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.308"
$ws.Range("E32").Value = "  +1.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.075"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05525"
$ws.Range("E34").Value = "  +6.39%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7324"
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.701"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  +1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.789"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4419"
$ws.Range("E40").Value = "  +0.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.33"
$ws.Range("E41").Value = "  +1.57%  "
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9989"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.895"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.90"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.585"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.715"
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "986.11"
$ws.Range("E49").Value = "  +10.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.058.81"
$ws.Range("E50").Value = "  +1.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.19"
$ws.Range("E51").Value = "  +0.82%  "
